$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Interior.ColorIndex = 6
Write-Output "done"
